# The report used to contain two parallel tables (A:B = raw per-row data,
# C:E = a "Jobs / Salary / Count" pivot-style summary). This update removes
# the now-unused C:E summary table, keeping only the Jobs Name / Salary (USD)
# table in columns A:B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old C:E summary table contents (rows 1-16 had data there).
# Clearing the cells (rather than just blanking values) drops the <c> nodes
# entirely, shrinking the sheet dimension down to A1:B20.
$ws.Range("C1:E16").ClearContents() | Out-Null

# Column A's header used to read "Jobs" (shared with the now-removed C1
# header); rename it to "Jobs Name" to match the remaining table's intent.
$ws.Range("A1").Value = "Jobs Name"

# Reflect the last on-screen selection over the (now cleared) old summary
# range, matching the saved view state.
$ws.Range("C1:E11").Select() | Out-Null
